$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 1, shifting all data down by one row.
$ws.Rows.Item(1).Insert()

# Populate the new header row with column titles.
$ws.Range("A1").Value = "idade"
$ws.Range("B1").Value = "volume"

# Match the final selection recorded in the workbook.
$ws.Range("B2").Select()
